# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# Add a new named cell style ("title_") -- bold + underlined variant of
# the existing "title" style -- mirroring the style table update that
# accompanied this revision.
$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true

# Remove the MSME breakdown row (Micro / SMEs / MSMEs) that used to live
# in row 5; only the Niger header and section title remain.
$ws.Rows.Item(5).Delete()
